$d = $word.ActiveDocument

# Locate the paragraph that contains the "Ver no Jupiter ..." text and
# work from there: the blurb consists of an empty paragraph, the
# "Ver no Jupiter Salvar em pdf Salvar em docx" paragraph, and the
# "(c) 2020 ... Creative Commons Attribution" paragraph, immediately
# followed by the LOQ4073 requirement line and a trailing blank
# paragraph before the page break. We need to drop the "Ver no
# Jupiter..." and copyright paragraphs together with one of the two
# blank paragraphs that surround them, leaving a single blank
# paragraph behind (matching the upstream diff).

$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Ver no Jupiter")) {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $copyrightPara = $target.Next()
    $trailingBlank = $copyrightPara.Next()

    # Remove the trailing blank paragraph that follows the copyright
    # line first (merging its mark forward keeps the untouched blank
    # paragraph that sits right after "LOQ4073..." intact).
    $trailingBlank.Range.Delete()
    $copyrightPara.Range.Delete()
    $target.Range.Delete()
}

Write-Output "done"
